$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 was stored as a text "2" in the original file -> convert to a real number
$ws.Range("B9").Value = 2

# Add new row 10 with the annotation data
$ws.Range("A10").Value = "Ying Tang"

# B10 holds the text "5" (not a number) - force text type, then clear the
# resulting formatting so no stray style is left behind on the cell
$b10 = $ws.Range("B10")
$b10.Value = "'5"
$b10.ClearFormats()

$ws.Range("C10").Value = "Thank, thoughtful "
$ws.Range("D10").Value = "ACK"
$ws.Range("E10").Value = "OTH"
$ws.Range("F10").Value = "41c93df3-3a59-4ce4-b94b-f420b7540586"
$ws.Range("G10").Value = "SJ19eUg0-_annotated.xlsx"
$ws.Range("H10").Value = "Thank the reviewer for the thoughtful feedback."
